$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset for rows 2..22 (columns D..I). Row 2 keeps the old row-31
# values; rows 3..22 carry a fresh block of simulation results that
# continues the trigger sweep past 3000.
$data = @(
        @(3000,348499,19,15,4,-24206),
        @(3100,341900,20,14,6,-25177),
        @(3200,325976,17,14,3,-23506),
        @(3300,297727,20,13,7,-22314),
        @(3400,293229,19,13,6,-21676),
        @(3500,283844,14,12,2,-20141),
        @(3800,282968,17,12,5,-20484),
        @(3600,282471,16,12,4,-21515),
        @(3700,276465,18,12,6,-20484),
        @(3900,271853,14,11,3,-19367),
        @(4000,270650,14,11,3,-19367),
        @(4100,266058,13,11,2,-19018),
        @(4200,256386,15,11,4,-18652),
        @(4300,252426,14,10,4,-18147),
        @(4400,247923,14,10,4,-16891),
        @(4500,245722,15,10,5,-17367),
        @(4700,234645,15,10,5,-17853),
        @(4600,232216,16,9,7,-17853),
        @(4800,229074,13,9,4,-16822),
        @(4900,222044,10,9,1,-14648),
        @(5000,219821,9,9,0,-14648)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $vals = $data[$idx]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
    $ws.Range("H$row").Value = $vals[4]
    $ws.Range("I$row").Value = $vals[5]
}

# The sweep now stops at row 22 - remove the old trailing rows 23..31.
$ws.Range("A23:I31").EntireRow.Delete()
